$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet 1: quality_comparison ---
# Give C1/D1 a distinct border style (top+bottom only for C1, top+bottom+right for D1)
$ws1.Range("C1").Borders.Item(7).LineStyle = 0
$ws1.Range("C1").Borders.Item(10).LineStyle = 0
$ws1.Range("D1").Borders.Item(7).LineStyle = 0

# Rename "fedcore" header to "approach"
$ws1.Range("C2").Value = "approach"

# --- Sheet 2: computational_comparison ---
$ws2.Range("C1").Borders.Item(7).LineStyle = 0
$ws2.Range("C1").Borders.Item(10).LineStyle = 0
$ws2.Range("D1").Borders.Item(7).LineStyle = 0

$ws2.Range("F1").Borders.Item(7).LineStyle = 0
$ws2.Range("F1").Borders.Item(10).LineStyle = 0
$ws2.Range("G1").Borders.Item(7).LineStyle = 0

$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# Remove the stray empty inline-string cell at G5
$ws2.Range("G5").ClearContents()
